$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.122445106506348
$ws.Range("B1").Value = 1.297770500183105
$ws.Range("C1").Value = 1.657012343406677
$ws.Range("D1").Value = 3.56758975982666
$ws.Range("E1").Value = 3.809174299240112
